$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Cade Cunningham"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Detroit Pistons"

$ws.Range("A3").Value = "Malik Monk"
$ws.Range("B3").Value = "PG,SG,SF"
$ws.Range("C3").Value = "Sacramento Kings"

$ws.Range("A4").Value = "LaMelo Ball"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Charlotte Hornets"

$ws.Range("A5").Value = "Brandon Miller"
$ws.Range("B5").Value = "SG,SF,PF"
$ws.Range("C5").Value = "Charlotte Hornets"

$ws.Range("A6").Value = "Jaden McDaniels"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Minnesota Timberwolves"

$ws.Range("A7").Value = "Ausar Thompson"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Detroit Pistons"

$ws.Range("A8").Value = "Julius Randle"
$ws.Range("B8").Value = "PF,C"
$ws.Range("C8").Value = "Minnesota Timberwolves"

$ws.Range("A9").Value = "Jonathan Isaac"
$ws.Range("B9").Value = "SF,PF"
$ws.Range("C9").Value = "Orlando Magic"

$ws.Range("A10").Value = "Anthony Davis"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "Los Angeles Lakers"

$ws.Range("A11").Value = "Bam Adebayo"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Miami Heat"

$ws.Range("A12").Value = "Cole Anthony"
$ws.Range("B12").Value = "PG"
$ws.Range("C12").Value = "Orlando Magic"

$ws.Range("A13").Value = "Cameron Johnson"
$ws.Range("B13").Value = "SF,PF"
$ws.Range("C13").Value = "Brooklyn Nets"

$ws.Range("A14").Value = "Damian Lillard"
$ws.Range("B14").Value = "PG"
$ws.Range("C14").Value = "Milwaukee Bucks"

$ws.Range("A15").Value = "Derrick White"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "Boston Celtics"

$ws.Range("A16").Value = "Naz Reid"
$ws.Range("B16").Value = "PF,C"
$ws.Range("C16").Value = "Minnesota Timberwolves"

$ws.Range("A17").Value = "Isaiah Hartenstein"
$ws.Range("B17").Value = "C"
$ws.Range("C17").Value = "Oklahoma City Thunder"

$ws.Range("A18").Value = "Jerami Grant"
$ws.Range("B18").Value = "SF,PF"
$ws.Range("C18").Value = "Portland Trail Blazers"

$ws.Range("A19").Value = "Collin Sexton"
$ws.Range("B19").Value = "PG,SG"
$ws.Range("C19").Value = "Utah Jazz"
